# Add 2022-Q3 data:
#  1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before "2022-Q2"),
#     populated with the new quarter's fund-holding data.
#  2. Insert a new row into the "总计" (summary) sheet for the "2022-Q3" entry,
#     shifting the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$newSheet.Name = "2022-Q3"

# Header row - bold, centered and bordered, matching the other quarter sheets
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108 # xlCenter
$newSheet.Range("B1:H1").VerticalAlignment = -4160   # xlTop
$newSheet.Range("B1:H1").Borders.LineStyle = 1

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row - B through G must stay text (matches "inlineStr" typing of the source data)
$newSheet.Range("B2:G2").NumberFormat = "@"

# Column A's running index uses the same bold/bordered look as the header
$newSheet.Range("A2").Font.Bold = $true
$newSheet.Range("A2").HorizontalAlignment = -4108 # xlCenter
$newSheet.Range("A2").VerticalAlignment = -4160   # xlTop
$newSheet.Range("A2").Borders.LineStyle = 1

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "513080"
$newSheet.Range("C2").Value = "华安法国CAC40ETF（QDII）"
$newSheet.Range("D2").Value = "0.58"
$newSheet.Range("E2").Value = "95.06"
$newSheet.Range("F2").Value = "6.68"
$newSheet.Range("G2").Value = "0.0387"
$newSheet.Range("H2").Value = 3

# ---------------------------------------------------------------------------
# 2) Insert a new row into "总计" for the "2022-Q3" figures, shifting the
#    remaining quarters down by one row.
# ---------------------------------------------------------------------------
$summary.Rows("2:2").Insert()

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.04
$summary.Cells.Item(2, 1).Value = 0

# The freshly-inserted row can pick up formatting from the header row above
# it; restore the plain data-row look (bold/bordered index cell in column A,
# unstyled B:D) by copying formats from an already-correct data row.
$summary.Range("A3:D3").Copy() | Out-Null
$summary.Range("A2:D2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-number the leading index column (A) sequentially for every data row,
# since inserting a row only shifts the old literal values rather than
# recomputing them.
for ($i = 0; $i -le 7; $i++) {
    $row = 2 + $i
    $summary.Cells.Item($row, 1).Value = $i
}
